$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '59.393.25'
$ws.Range('E2').Value = '  +3.98%  '

# Row 3
$ws.Range('D3').Value = '3.311.28'
$ws.Range('E3').Value = '  +1.12%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '406.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.94%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.581'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.629'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.12%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.01%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0980'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.52%  '

# Row 12
$ws.Range('E12').Value = '  +1.16%  '

# Row 13
$ws.Range('D13').Value = '3.829.15'
$ws.Range('E13').Value = '  +1.18%  '

# Row 14
$ws.Range('E14').Value = '  +4.43%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.21%  '

# Row 16
$ws.Range('D16').Value = '3.311.03'
$ws.Range('E16').Value = '  +1.36%  '

# Row 17
$ws.Range('E17').Value = '  -1.68%  '

# Row 18
$ws.Range('D18').Value = '59.090.88'
$ws.Range('E18').Value = '  +3.89%  '

# Row 19
$ws.Range('E19').Value = '  -2.88%  '

# Row 20
$ws.Range('E20').Value = '  -0.12%  '

# Row 21
$ws.Range('E21').Value = '  +5.62%  '

# Row 22
$ws.Range('E22').Value = '  -0.07%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '306.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.16%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.00%  '

# Row 25
$ws.Range('E25').Value = '  -0.83%  '

# Row 26
$ws.Range('E26').Value = '  +2.40%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.181'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.81%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.77'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.49%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.40'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.88%  '

# Row 31
$ws.Range('E31').Value = '  +2.68%  '

# Row 32
$ws.Range('E32').Value = '  +0.03%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.99%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '39.66'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.27%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0508'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.79%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.82%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.15'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.46%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.03%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.24%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '138.66'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.38%  '

# Row 42
$ws.Range('E42').Value = '  +2.01%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.91%  '

# Row 44
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.37%  '

# Row 45
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.81'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.13%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.280'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.38%  '

# Row 47
$ws.Range('E47').Value = '  +8.57%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.16%  '

# Row 49
$ws.Range('D49').Value = '2.200.97'
$ws.Range('E49').Value = '  +1.79%  '

# Row 50
$ws.Range('E50').Value = '  +0.13%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.22%  '
